$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.83119999999999
$ws.Range("E8").Value = 13.63409999999999
$ws.Range("C12").Value = -14.82640000000003
$ws.Range("E12").Value = 12.43739999999999
$ws.Range("E14").Value = 13.7269
$ws.Range("E22").Value = 11.97019999999999
